$wb = $excel.ActiveWorkbook

# ---- Sheet 1: treatment ----
$ws1 = $wb.Worksheets.Item("treatment")

# Updated values (rows 3-13) reflecting re-run meta-analysis estimates
$ws1.Range("B3").Value = 42.54213354518064
$ws1.Range("C3").Value = 64.98360952142842
$ws1.Range("D3").Value = 95.16011944566623
$ws1.Range("B4").Value = 0.3134547872257854
$ws1.Range("C4").Value = 0.8059130228725063
$ws1.Range("D4").Value = 1.526837100606182
$ws1.Range("B5").Value = 0.5703668455414779
$ws1.Range("C5").Value = 0.9145572694398032
$ws1.Range("D5").Value = 1.258818658355301
$ws1.Range("B7").Value = 20.22111161544328
$ws1.Range("C7").Value = 25.53905710956621
$ws1.Range("D7").Value = 31.82574352472265
$ws1.Range("B8").Value = 0.111712121023525
$ws1.Range("C8").Value = 0.2542928722656458
$ws1.Range("D8").Value = 0.4545099690195018
$ws1.Range("B9").Value = 0.3392407178687752
$ws1.Range("C9").Value = 0.5118290343848231
$ws1.Range("D9").Value = 0.6842732007438587
$ws1.Range("B11").Value = 7.548231660743239
$ws1.Range("C11").Value = 10.81097546473494
$ws1.Range("D11").Value = 15.00461574004915
$ws1.Range("B12").Value = 0.2285994735526786
$ws1.Range("C12").Value = 0.5712573783272646
$ws1.Range("D12").Value = 1.068105752310871
$ws1.Range("B13").Value = 0.486654228821725
$ws1.Range("C13").Value = 0.7693054335271695
$ws1.Range("D13").Value = 1.051937689563459

# New "Speed meta analysis" section (rows 14-17)
$ws1.Range("A14").Value = "Speed meta analysis"
$ws1.Range("A15").Value = "mean (km/day)"
$ws1.Range("A16").Value = "CoV² (RVAR)"
$ws1.Range("A17").Value = "CoV  (RSTD)"
$ws1.Range("B15").Value = 3.625695495129722
$ws1.Range("C15").Value = 4.367383454604379
$ws1.Range("D15").Value = 5.212968018226442
$ws1.Range("B16").Value = 0.04546823585390327
$ws1.Range("C16").Value = 0.1200519960437265
$ws1.Range("D16").Value = 0.2302092668221222
$ws1.Range("B17").Value = 0.2174120584257805
$ws1.Range("C17").Value = 0.3532761090516405
$ws1.Range("D17").Value = 0.489205100840114

# ---- Sheet 2: control ----
$ws2 = $wb.Worksheets.Item("control")

# Updated values (rows 3-13) reflecting re-run meta-analysis estimates
$ws2.Range("B3").Value = 32.79463319010795
$ws2.Range("C3").Value = 42.00147592616612
$ws2.Range("D3").Value = 52.98013858976994
$ws2.Range("B4").Value = 0.1210038032932299
$ws2.Range("C4").Value = 0.2799244435321588
$ws2.Range("D4").Value = 0.5043441193402219
$ws2.Range("B5").Value = 0.3532361919942047
$ws2.Range("C5").Value = 0.537261964299384
$ws2.Range("D5").Value = 0.7211559058305435
$ws2.Range("B7").Value = 18.20460511584151
$ws2.Range("C7").Value = 22.23121098640574
$ws2.Range("D7").Value = 26.87987521833361
$ws2.Range("B8").Value = 0.08327953093920003
$ws2.Range("C8").Value = 0.186876953668574
$ws2.Range("D8").Value = 0.3316250264200705
$ws2.Range("B9").Value = 0.2927818026158576
$ws2.Range("C9").Value = 0.4385840343338146
$ws2.Range("D9").Value = 0.5842498289817237
$ws2.Range("B11").Value = 5.051322294049788
$ws2.Range("C11").Value = 5.909450055679693
$ws2.Range("D11").Value = 6.867787399791214
$ws2.Range("B12").Value = 0.04522646809648333
$ws2.Range("C12").Value = 0.1072663177646837
$ws2.Range("D12").Value = 0.1956343536194577
$ws2.Range("B13").Value = 0.2161158490283523
$ws2.Range("C13").Value = 0.3328296182750619
$ws2.Range("D13").Value = 0.4494825546923926

# New "Speed meta analysis" section (rows 14-17)
$ws2.Range("A14").Value = "Speed meta analysis"
$ws2.Range("A15").Value = "mean (km/day)"
$ws2.Range("A16").Value = "CoV² (RVAR)"
$ws2.Range("A17").Value = "CoV  (RSTD)"
$ws2.Range("B15").Value = 3.362856235691671
$ws2.Range("C15").Value = 3.806741249265726
$ws2.Range("D15").Value = 4.290008097159491
$ws2.Range("B16").Value = 0.02098895657273051
$ws2.Range("C16").Value = 0.0578806186685388
$ws2.Range("D16").Value = 0.1131540085843296
$ws2.Range("B17").Value = 0.1479196194664224
$ws2.Range("C17").Value = 0.2456387839265089
$ws2.Range("D17").Value = 0.3434514295130752
